# Gaussian Quadrature Scheme export: rename sheet, touch a few fp values,
# and append a new averaged-intensity row (row 16) to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name no longer carries the raw .xpc suffix)
$ws.Name = "BetaFiberA"

# Minor floating point refinements in row 13 (re-saved from the upstream calc)
$ws.Range("F13").Value = 0.99096221622741
$ws.Range("K13").Value = 0.9926970675175544
$ws.Range("O13").Value = 0.9934342321779698

# Append the new "HexGrid-60degTilt5degRes" averaged row (row 16)
# Copy row 15's label cells first so formatting/shared-string reuse matches,
# then overwrite with the new row's own values.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14

$ws.Range("B15").Copy($ws.Range("B16"))

$ws.Range("C16").Value = 1.210448842380348
$ws.Range("D16").Value = 0.9085960483277356
$ws.Range("E16").Value = 1.05656414440336
$ws.Range("F16").Value = 0.9104558492700081
$ws.Range("G16").Value = 1.210448842380348
$ws.Range("H16").Value = 0.9085960483277356
$ws.Range("I16").Value = 1.067917046248743
$ws.Range("J16").Value = 0.91967140095761
$ws.Range("K16").Value = 1.021511881579576
$ws.Range("L16").Value = 0.8813047786622972
$ws.Range("M16").Value = 1.210448842380348
$ws.Range("N16").Value = 0.9825800963655478
$ws.Range("O16").Value = 1.021516221095363
$ws.Range("P16").Value = 0.9970587489787096
